# Add season record columns (Wins, Losses, Ties) to the DET_1992 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers should look like the other header
# cells (bold, centered, bordered) so copy the format from the existing
# header cell AC1 before writing the new labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-36): every player on this roster shares the same team
# season record: 75 wins, 87 losses, 0 ties.
$wins = 75
$losses = 87
$ties = 0

for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
